$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.803.26"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "3.296.46"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "185.92"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").Value = "581.14"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").Value = "3.877.27"
$ws.Range("E12").Value = "  +1.29%  "

$ws.Range("E13").Value = "  -2.05%  "

$ws.Range("D14").Value = "27.41"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").Value = "67.993.38"
$ws.Range("E15").Value = "  +0.34%  "

$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").Value = "3.301.85"
$ws.Range("E17").Value = "  +1.46%  "

$ws.Range("D18").Value = "447.43"
$ws.Range("E18").Value = "  +12.20%  "

$ws.Range("D19").Value = "5.70"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "13.50"
$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("D21").Value = "7.71"
$ws.Range("E21").Value = "  +2.02%  "

$ws.Range("D22").Value = "74.86"
$ws.Range("E22").Value = "  +5.56%  "

$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "3.466.21"
$ws.Range("E24").Value = "  +1.76%  "

$ws.Range("D25").Value = "0.513"
$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("E26").Value = "  +1.14%  "

$ws.Range("E27").Value = "  +0.42%  "

$ws.Range("D28").Value = "9.09"
$ws.Range("E28").Value = "  -4.30%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("D30").Value = "1.97"
$ws.Range("E30").Value = "  +1.27%  "

$ws.Range("D31").Value = "22.83"
$ws.Range("E31").Value = "  +0.93%  "

$ws.Range("E32").Value = "  -2.06%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").Value = "6.79"
$ws.Range("E35").Value = "  -1.95%  "

$ws.Range("E36").Value = "  +4.36%  "

$ws.Range("D37").Value = "163.50"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("D38").Value = "1.86"
$ws.Range("E38").Value = "  -1.41%  "

$ws.Range("D39").Value = "26.97"
$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("D40").Value = "4.50"
$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("E41").Value = "  -2.78%  "

$ws.Range("D42").Value = "6.40"
$ws.Range("E42").Value = "  +2.21%  "

$ws.Range("D43").Value = "2.697.84"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("D44").Value = "40.70"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("D45").Value = "0.0672"
$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("D46").Value = "2.40"

$ws.Range("D47").Value = "24.64"
$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("D48").Value = "325.14"
$ws.Range("E48").Value = "  -3.04%  "

$ws.Range("D49").Value = "0.0275"
$ws.Range("E49").Value = "  +0.78%  "

$ws.Range("D50").Value = "31.48"
$ws.Range("E50").Value = "  +3.04%  "

$ws.Range("D51").Value = "0.987"
$ws.Range("E51").Value = "  +2.00%  "
